# Weekly cryptos data refresh (GitHub Actions bot)
# Updates Price (D) and Volume(1h) (E) columns; rows 42/43 and 46/47
# also swap rank position (B/C/D/E) between two coins.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.249.16'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.34%  '

$ws.Range("D3").Value = '''1.605.68'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.10%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '''212.73'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.04%  '

$ws.Range("E6").Value = '  -0.11%  '

$ws.Range("D7").Value = '''0.487'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("E8").Value = '  +0.65%  '

$ws.Range("D9").Value = '''0.0614'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.27%  '

$ws.Range("D10").Value = '''18.40'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.18%  '

$ws.Range("E11").Value = '  -0.60%  '

$ws.Range("D12").Value = '''1.828.14'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.05%  '

$ws.Range("D13").Value = '''1.613.99'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.61%  '

$ws.Range("E14").Value = '  +0.51%  '

$ws.Range("E15").Value = '  +0.61%  '

$ws.Range("D16").Value = '''26.218.39'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.25%  '

$ws.Range("D17").Value = '''62.05'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.58%  '

$ws.Range("E18").Value = '  +0.85%  '

$ws.Range("D20").Value = '''200.26'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.02%  '

$ws.Range("D21").Value = '''4.27'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.66%  '

$ws.Range("D22").Value = '''9.32'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.16%  '

$ws.Range("D24").Value = '''1.88'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.89%  '

$ws.Range("D25").Value = '''144.22'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.91%  '

$ws.Range("E26").Value = '  -0.09%  '

$ws.Range("E27").Value = '  -2.17%  '

$ws.Range("D28").Value = '''15.19'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.01%  '

$ws.Range("D29").Value = '''6.57'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.04%  '

$ws.Range("E30").Value = '  +4.49%  '

$ws.Range("D31").Value = '''1.18'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.68%  '

$ws.Range("E32").Value = '  +2.73%  '

$ws.Range("E33").Value = '  -1.23%  '

$ws.Range("E34").Value = '  +1.12%  '

$ws.Range("E35").Value = '  +1.05%  '

$ws.Range("D36").Value = '''1.166.54'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +4.23%  '

$ws.Range("E37").Value = '  +3.40%  '

$ws.Range("E38").Value = '  -0.11%  '

$ws.Range("E39").Value = '  +0.84%  '

$ws.Range("D40").Value = '''0.785'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.25%  '

$ws.Range("D41").Value = '''0.498'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.21%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '''0.785'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.53%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''5.34'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +4.05%  '

$ws.Range("D44").Value = '''1.739.60'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.01%  '

$ws.Range("D45").Value = '''92.15'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.81%  '

$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").Value = '''1.54'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.56%  '

$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '''0.0₆0105'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +13.60%  '

$ws.Range("D48").Value = '''54.05'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.04%  '

$ws.Range("D49").Value = '''0.0505'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.21%  '

$ws.Range("E50").Value = '  -0.41%  '

$ws.Range("E51").Value = '  -0.18%  '
